$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RPM sample input (speed control tune data) and let the
# rad/s formula in B4 recalculate off the new value.
$ws.Range("A4").Value = 120

# New encoder debug row: a literal "<" marker plus an unrelated scratch
# ratio calc parked out in column F.
$ws.Range("A5").Value = "<"
$ws.Range("F5").Formula = "=12/5.6"

# Extra tune-data ratio calc further down the sheet.
$ws.Range("A10").Formula = "=57/65"

# Leave the selection where the new data entry ended up.
$ws.Range("A5").Select() | Out-Null
